$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) for rows 16-37 currently lists the period
# codes in descending order (2003 down to 1806). This update re-sorts them
# in ascending order (1806 up to 2003) as part of refreshing the account
# statement database.
$periods = @("1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}
